$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "40"

$t.Cell(6,1).Range.Text = "0.00052"
$t.Cell(7,1).Range.Text = "0.00024"
$t.Cell(8,1).Range.Text = "0.00005"
$t.Cell(9,1).Range.Text = "0.00044"
$t.Cell(10,1).Range.Text = "0.00049"
$t.Cell(11,1).Range.Text = "0.00050"
$t.Cell(12,1).Range.Text = "0.00959"

$t.Cell(44,1).Range.Text = "99.95"
$t.Cell(45,1).Range.Text = "0.01"
$t.Cell(46,1).Range.Text = "19"
